$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")

# Switch TS_001 "Open the Browser" step (row 2) to use headless Chrome
# instead of Mozilla as its Data Set value.
$ws.Range("H2").Value = "Chrome"

# Remove the obsolete waitFor5 test step (TS_032, row 32)
$ws.Rows.Item(32).Delete()

$ws.Range("H2").Select()
